# fix login user upload, add role user upload
# Adds a "Role" column (E) with a "company viewer" default for each data
# row, backed by an in-cell dropdown (data validation list) offering the
# three known roles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("E1").Value = "Role"

# Default role value for each existing data row
$ws.Range("E2").Value = "company viewer"
$ws.Range("E3").Value = "company viewer"
$ws.Range("E4").Value = "company viewer"

# Give the new column a comfortable width (matches the other bespoke
# column widths already present in the sheet)
$ws.Columns.Item(5).ColumnWidth = 18.33

# Dropdown list validation on the Role cells
$roleRange = $ws.Range("E2:E4")
$roleRange.Validation.Add(3, 1, 1, '"company coordinator, company as, company viewer"') | Out-Null

# Leave the cursor where the user last clicked while filling this in
$ws.Range("I17").Select() | Out-Null
